$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record progress dates for Exercises 7-12 (rows 10-17, column B). These use
# the same date (2018-01-10, serial 43110) and the same date format/centered
# alignment that the earlier rows (3-9, serial 43109) already use.
$rng = $ws.Range("B10:B17")
$rng.Value = 43110
$rng.NumberFormat = "d-mmm"
$rng.HorizontalAlignment = -4108

# Move the active cell/selection to the last entry just filled in (B17)
$ws.Range("B17").Select()

# Best-effort: nudge the saved window position (purely cosmetic screen
# position; harmless no-op if the host doesn't expose it)
try { $excel.ActiveWindow.Left = 12700 } catch {}
